$wb = $excel.ActiveWorkbook

# Update "想去人数" (attendance count) figures on both the "展览" sheet
# and the "全部类型" sheet, which carry duplicate rows for the same events.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4660
    $ws.Range("F3").Value = 137
    $ws.Range("F4").Value = 809
}
